$wb = $excel.ActiveWorkbook

# ---- Sheet 1: Overview ----
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Hyperlinks.Delete()

$ws1.Range("B2").Value = 'Handed back: in sync with en-US'
$ws1.Range("C2").Value = 'Handed back: in sync with en-US'
$ws1.Range("D2").Value = '2016-03-22 06:37:29'

$ws1.Range("B3").Value = 'Handed back: in sync with en-US'
$ws1.Range("C3").Value = 'Handed back: in sync with en-US'
$ws1.Range("D3").Value = '2016-03-22 06:35:35'

$ws1.Range("B4").Value = 'Handed back: in sync with en-US'
$ws1.Range("C4").Value = 'Handed back: in sync with en-US'
$ws1.Range("D4").Value = '2016-03-22 06:37:29'

$ws1.Hyperlinks.Add($ws1.Range("A2"), 'https://github.com/OpenLocalizationTest/oltest/blob/a3204d45895904f54c6793dd67fc76f847b74c46/e2e/2e2f5251-46c3-4d95-998c-0c30dcc29ee7.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '2e2f5251-46c3-4d95-998c-0c30dcc29ee7.md')
$ws1.Hyperlinks.Add($ws1.Range("A3"), 'https://github.com/OpenLocalizationTest/oltest/blob/4dfc8a3422e6f711e2b497c394ca86cd16d094cf/e2e/4a88754c-d4d3-40bf-883e-05388a36cbfc.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '4a88754c-d4d3-40bf-883e-05388a36cbfc.md')
$ws1.Hyperlinks.Add($ws1.Range("A4"), 'https://github.com/OpenLocalizationTest/oltest/blob/a3204d45895904f54c6793dd67fc76f847b74c46/e2e/57569c52-c9e2-42eb-a531-e2dcac98a1a6.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '57569c52-c9e2-42eb-a531-e2dcac98a1a6.md')

# ---- Sheet: zh-cn ----
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Delete()

$ws2.Range("B2").Value = '.md'
$ws2.Range("C2").Value = 'Handed back: in sync with en-US'
$ws2.Range("E2").Value = '2016-03-22 06:37:21'
$ws2.Range("H2").Value = '2016-03-22 06:38:03'
$ws2.Range("J2").Value = 'Include'

$ws2.Range("B3").Value = '.md'
$ws2.Range("C3").Value = 'Handed back: in sync with en-US'
$ws2.Range("E3").Value = '2016-03-22 06:35:26'
$ws2.Range("H3").Value = '2016-03-22 06:36:16'
$ws2.Range("J3").Value = 'Include'

$ws2.Range("B4").Value = '.md'
$ws2.Range("C4").Value = 'Handed back: in sync with en-US'
$ws2.Range("E4").Value = '2016-03-22 06:37:21'
$ws2.Range("H4").Value = '2016-03-22 06:38:03'
$ws2.Range("J4").Value = 'Include'

$ws2.Hyperlinks.Add($ws2.Range("A2"), 'https://github.com/OpenLocalizationTest/oltest/blob/a3204d45895904f54c6793dd67fc76f847b74c46/e2e/2e2f5251-46c3-4d95-998c-0c30dcc29ee7.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '2e2f5251-46c3-4d95-998c-0c30dcc29ee7.md')
$ws2.Hyperlinks.Add($ws2.Range("D2"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a39dc479226e2904947c0239fd379d20ed352052/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/2e2f5251-46c3-4d95-998c-0c30dcc29ee7.d3b14ffdec5fa8dc65838d1d7121fdbe518b5765.zh-cn.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '2e2f5251-46c3-4d95-998c-0c30dcc29ee7.d3b14ffdec5fa8dc65838d1d7121fdbe518b5765.zh-cn.xlf')
$ws2.Hyperlinks.Add($ws2.Range("F2"), 'https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/148bbca0ec6c1de92514e5f43fc90f2bfe055635/e2e/2e2f5251-46c3-4d95-998c-0c30dcc29ee7.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '2e2f5251-46c3-4d95-998c-0c30dcc29ee7.md')
$ws2.Hyperlinks.Add($ws2.Range("G2"), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/f511a18c27a3e117220f8beb1b37e75ae6155261/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/2e2f5251-46c3-4d95-998c-0c30dcc29ee7.d3b14ffdec5fa8dc65838d1d7121fdbe518b5765.zh-cn.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '2e2f5251-46c3-4d95-998c-0c30dcc29ee7.d3b14ffdec5fa8dc65838d1d7121fdbe518b5765.zh-cn.xlf')

$ws2.Hyperlinks.Add($ws2.Range("A3"), 'https://github.com/OpenLocalizationTest/oltest/blob/4dfc8a3422e6f711e2b497c394ca86cd16d094cf/e2e/4a88754c-d4d3-40bf-883e-05388a36cbfc.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '4a88754c-d4d3-40bf-883e-05388a36cbfc.md')
$ws2.Hyperlinks.Add($ws2.Range("D3"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a39dc479226e2904947c0239fd379d20ed352052/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/4a88754c-d4d3-40bf-883e-05388a36cbfc.5e8169c68a7a747323cf9ae25f88cdd24a801999.zh-cn.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '4a88754c-d4d3-40bf-883e-05388a36cbfc.5e8169c68a7a747323cf9ae25f88cdd24a801999.zh-cn.xlf')
$ws2.Hyperlinks.Add($ws2.Range("F3"), 'https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/148bbca0ec6c1de92514e5f43fc90f2bfe055635/e2e/4a88754c-d4d3-40bf-883e-05388a36cbfc.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '4a88754c-d4d3-40bf-883e-05388a36cbfc.md')
$ws2.Hyperlinks.Add($ws2.Range("G3"), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/f511a18c27a3e117220f8beb1b37e75ae6155261/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/4a88754c-d4d3-40bf-883e-05388a36cbfc.5e8169c68a7a747323cf9ae25f88cdd24a801999.zh-cn.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '4a88754c-d4d3-40bf-883e-05388a36cbfc.5e8169c68a7a747323cf9ae25f88cdd24a801999.zh-cn.xlf')

$ws2.Hyperlinks.Add($ws2.Range("A4"), 'https://github.com/OpenLocalizationTest/oltest/blob/a3204d45895904f54c6793dd67fc76f847b74c46/e2e/57569c52-c9e2-42eb-a531-e2dcac98a1a6.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '57569c52-c9e2-42eb-a531-e2dcac98a1a6.md')
$ws2.Hyperlinks.Add($ws2.Range("D4"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a39dc479226e2904947c0239fd379d20ed352052/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/2e2f5251-46c3-4d95-998c-0c30dcc29ee7.d3b14ffdec5fa8dc65838d1d7121fdbe518b5765.zh-cn.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '2e2f5251-46c3-4d95-998c-0c30dcc29ee7.d3b14ffdec5fa8dc65838d1d7121fdbe518b5765.zh-cn.xlf')
$ws2.Hyperlinks.Add($ws2.Range("F4"), 'https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/148bbca0ec6c1de92514e5f43fc90f2bfe055635/e2e/2e2f5251-46c3-4d95-998c-0c30dcc29ee7.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '2e2f5251-46c3-4d95-998c-0c30dcc29ee7.md')
$ws2.Hyperlinks.Add($ws2.Range("G4"), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/f511a18c27a3e117220f8beb1b37e75ae6155261/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/2e2f5251-46c3-4d95-998c-0c30dcc29ee7.d3b14ffdec5fa8dc65838d1d7121fdbe518b5765.zh-cn.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '2e2f5251-46c3-4d95-998c-0c30dcc29ee7.d3b14ffdec5fa8dc65838d1d7121fdbe518b5765.zh-cn.xlf')

# ---- Sheet: de-de ----
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Hyperlinks.Delete()

$ws3.Range("B2").Value = '.md'
$ws3.Range("C2").Value = 'Handed back: in sync with en-US'
$ws3.Range("E2").Value = '2016-03-22 06:37:29'
$ws3.Range("H2").Value = '2016-03-22 06:38:17'
$ws3.Range("J2").Value = 'Include'

$ws3.Range("B3").Value = '.md'
$ws3.Range("C3").Value = 'Handed back: in sync with en-US'
$ws3.Range("E3").Value = '2016-03-22 06:35:35'
$ws3.Range("H3").Value = '2016-03-22 06:36:30'
$ws3.Range("J3").Value = 'Include'

$ws3.Range("B4").Value = '.md'
$ws3.Range("C4").Value = 'Handed back: in sync with en-US'
$ws3.Range("E4").Value = '2016-03-22 06:37:29'
$ws3.Range("H4").Value = '2016-03-22 06:38:17'
$ws3.Range("J4").Value = 'Include'

$ws3.Hyperlinks.Add($ws3.Range("A2"), 'https://github.com/OpenLocalizationTest/oltest/blob/a3204d45895904f54c6793dd67fc76f847b74c46/e2e/2e2f5251-46c3-4d95-998c-0c30dcc29ee7.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '2e2f5251-46c3-4d95-998c-0c30dcc29ee7.md')
$ws3.Hyperlinks.Add($ws3.Range("D2"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a39dc479226e2904947c0239fd379d20ed352052/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/2e2f5251-46c3-4d95-998c-0c30dcc29ee7.d3b14ffdec5fa8dc65838d1d7121fdbe518b5765.de-de.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '2e2f5251-46c3-4d95-998c-0c30dcc29ee7.d3b14ffdec5fa8dc65838d1d7121fdbe518b5765.de-de.xlf')
$ws3.Hyperlinks.Add($ws3.Range("F2"), 'https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/148bbca0ec6c1de92514e5f43fc90f2bfe055635/e2e/2e2f5251-46c3-4d95-998c-0c30dcc29ee7.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '2e2f5251-46c3-4d95-998c-0c30dcc29ee7.md')
$ws3.Hyperlinks.Add($ws3.Range("G2"), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/f511a18c27a3e117220f8beb1b37e75ae6155261/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/2e2f5251-46c3-4d95-998c-0c30dcc29ee7.d3b14ffdec5fa8dc65838d1d7121fdbe518b5765.de-de.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '2e2f5251-46c3-4d95-998c-0c30dcc29ee7.d3b14ffdec5fa8dc65838d1d7121fdbe518b5765.de-de.xlf')

$ws3.Hyperlinks.Add($ws3.Range("A3"), 'https://github.com/OpenLocalizationTest/oltest/blob/4dfc8a3422e6f711e2b497c394ca86cd16d094cf/e2e/4a88754c-d4d3-40bf-883e-05388a36cbfc.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '4a88754c-d4d3-40bf-883e-05388a36cbfc.md')
$ws3.Hyperlinks.Add($ws3.Range("D3"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a39dc479226e2904947c0239fd379d20ed352052/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/4a88754c-d4d3-40bf-883e-05388a36cbfc.5e8169c68a7a747323cf9ae25f88cdd24a801999.de-de.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '4a88754c-d4d3-40bf-883e-05388a36cbfc.5e8169c68a7a747323cf9ae25f88cdd24a801999.de-de.xlf')
$ws3.Hyperlinks.Add($ws3.Range("F3"), 'https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/148bbca0ec6c1de92514e5f43fc90f2bfe055635/e2e/4a88754c-d4d3-40bf-883e-05388a36cbfc.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '4a88754c-d4d3-40bf-883e-05388a36cbfc.md')
$ws3.Hyperlinks.Add($ws3.Range("G3"), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/f511a18c27a3e117220f8beb1b37e75ae6155261/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/4a88754c-d4d3-40bf-883e-05388a36cbfc.5e8169c68a7a747323cf9ae25f88cdd24a801999.de-de.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '4a88754c-d4d3-40bf-883e-05388a36cbfc.5e8169c68a7a747323cf9ae25f88cdd24a801999.de-de.xlf')

$ws3.Hyperlinks.Add($ws3.Range("A4"), 'https://github.com/OpenLocalizationTest/oltest/blob/a3204d45895904f54c6793dd67fc76f847b74c46/e2e/57569c52-c9e2-42eb-a531-e2dcac98a1a6.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '57569c52-c9e2-42eb-a531-e2dcac98a1a6.md')
$ws3.Hyperlinks.Add($ws3.Range("D4"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a39dc479226e2904947c0239fd379d20ed352052/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/2e2f5251-46c3-4d95-998c-0c30dcc29ee7.d3b14ffdec5fa8dc65838d1d7121fdbe518b5765.de-de.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '2e2f5251-46c3-4d95-998c-0c30dcc29ee7.d3b14ffdec5fa8dc65838d1d7121fdbe518b5765.de-de.xlf')
$ws3.Hyperlinks.Add($ws3.Range("F4"), 'https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/148bbca0ec6c1de92514e5f43fc90f2bfe055635/e2e/2e2f5251-46c3-4d95-998c-0c30dcc29ee7.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '2e2f5251-46c3-4d95-998c-0c30dcc29ee7.md')
$ws3.Hyperlinks.Add($ws3.Range("G4"), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/f511a18c27a3e117220f8beb1b37e75ae6155261/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/2e2f5251-46c3-4d95-998c-0c30dcc29ee7.d3b14ffdec5fa8dc65838d1d7121fdbe518b5765.de-de.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '2e2f5251-46c3-4d95-998c-0c30dcc29ee7.d3b14ffdec5fa8dc65838d1d7121fdbe518b5765.de-de.xlf')

